$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor has signed off the timesheet: fill in the supervisor's name
# next to "Supervisor Name:" (row 6).
$ws.Range("G6").Value = "Ankita Gangotra"

# Fill in the "Supervisor Signature" block (row 27): initials + sign-off date.
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").Value = "6/10/2014"
$ws.Range("D27").NumberFormat = "mm-dd-yy"

# Match the author's final selection (the supervisor date field).
$ws.Range("D27:E27").Select()
